$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.55"
$ws.Range("E2").Value = "'0.16%"
$ws.Range("E3").Value = "'-0.50%"
$ws.Range("D4").Value = "'5.037"
$ws.Range("E4").Value = "'-1.27%"
$ws.Range("D5").Value = "'0.08036"
$ws.Range("E5").Value = "'-0.59%"
$ws.Range("D6").Value = "'1.890"
$ws.Range("E6").Value = "'-2.41%"
$ws.Range("E7").Value = "'-0.89%"
$ws.Range("D8").Value = "'7.777"
$ws.Range("E8").Value = "'0.24%"
$ws.Range("D9").Value = "'0.9207"
$ws.Range("E9").Value = "'-0.70%"
$ws.Range("D10").Value = "'0.1266"
$ws.Range("E10").Value = "'-6.26%"
$ws.Range("D11").Value = "'0.1911"
$ws.Range("E11").Value = "'0.23%"
$ws.Range("D12").Value = "'0.09120"
$ws.Range("E12").Value = "'-1.05%"
$ws.Range("D13").Value = "'0.03465"
$ws.Range("E13").Value = "'1.69%"
$ws.Range("D14").Value = "'0.09855"
$ws.Range("E14").Value = "'0.31%"
$ws.Range("D15").Value = "'0.001415"
$ws.Range("E15").Value = "'0.78%"
$ws.Range("D16").Value = "'0.006240"
$ws.Range("E16").Value = "'5.68%"
$ws.Range("D17").Value = "'3.815"
$ws.Range("E17").Value = "'7.35%"
$ws.Range("D18").Value = "'3.348"
$ws.Range("E18").Value = "'12.86%"
$ws.Range("D19").Value = "'0.3418"
$ws.Range("E20").Value = "'1.16%"
$ws.Range("D21").Value = "'5.169"
$ws.Range("E21").Value = "'5.48%"
$ws.Range("E22").Value = "'-11.41%"
$ws.Range("D23").Value = "'0.04417"
$ws.Range("E23").Value = "'0.28%"
$ws.Range("E24").Value = "'1.16%"
$ws.Range("D25").Value = "'0.004614"
$ws.Range("E25").Value = "'-4.01%"
$ws.Range("D27").Value = "'0.0001252"
$ws.Range("E27").Value = "'-3.81%"
$ws.Range("E28").Value = "'42.04%"
$ws.Range("D39").Value = "'0.01947"
$ws.Range("E39").Value = "'-3.07%"
$ws.Range("D40").Value = "'0.05388"
$ws.Range("E40").Value = "'9.72%"
$ws.Range("D41").Value = "'0.007604"
$ws.Range("E41").Value = "'-0.46%"
$ws.Range("D42").Value = "'0.01018"
$ws.Range("E42").Value = "'-0.65%"
$ws.Range("D43").Value = "'0.1351"
$ws.Range("E43").Value = "'-1.76%"
$ws.Range("D44").Value = "'0.002154"
$ws.Range("E44").Value = "'2.41%"
$ws.Range("D45").Value = "'0.009650"
$ws.Range("E45").Value = "'-11.38%"
$ws.Range("D46").Value = "'0.00006121"
$ws.Range("E46").Value = "'-4.14%"
$ws.Range("E47").Value = "'-0.03%"
$ws.Range("D48").Value = "'65.22"
$ws.Range("E48").Value = "'0.85%"
$ws.Range("E49").Value = "'39.29%"
$ws.Range("E50").Value = "'-0.03%"
$ws.Range("E51").Value = "'-0.03%"
